$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Special")

# --- B11: "1. Create Task: ..." -----------------------------------------
# Replace the last (red) bullet so it reads the new sentence, and drop the
# red highlight color from the whole body (it was selected + "No Fill"
# colored back to automatic in the real edit).
$cell11 = $ws.Range("B11")

$t11_1 = "1. Create Task: "
$t11_2 = "`n- Thiết lập trạng thái mặc định là NOT_STARTED nếu không có`n- Gửi thông báo tới project manager, project member về việc tạo công việc mới`n- Thời lượng Task đang k đc quá 1 năm"
$t11_3 = "`n"
$t11_4 = "- Giới hạn thời lượng Task phải thuộc thời lượng Project"

$full11 = $t11_1 + $t11_2 + $t11_3 + $t11_4
$cell11.Value = $full11

$p11_1 = 1
$p11_2 = $p11_1 + $t11_1.Length
$p11_3 = $p11_2 + $t11_2.Length
$p11_4 = $p11_3 + $t11_3.Length

$cell11.Characters($p11_2, $t11_2.Length).Font.ColorIndex = -4105
$cell11.Characters($p11_3, $t11_3.Length).Font.Bold = $true
$cell11.Characters($p11_3, $t11_3.Length).Font.ColorIndex = -4105
$cell11.Characters($p11_4, $t11_4.Length).Font.ColorIndex = -4105

# --- B17: "1. Create subtask ..." ---------------------------------------
$cell17 = $ws.Range("B17")

$t17_1 = "1. Create subtask`n"
$t17_2 = "- Gửi thông báo cho người được giao việc phụ`n- Giới hạn thời lượng Subtask phải thuộc thời lượng Task"

$full17 = $t17_1 + $t17_2
$cell17.Value = $full17

$p17_1 = 1
$p17_2 = $p17_1 + $t17_1.Length

$cell17.Characters($p17_2, $t17_2.Length).Font.ColorIndex = -4105
